# Quarterly indexing esoteric bug-fix operation
# The dates stored in column A (rows 2-73) represent the start of a
# reporting quarter, but they need to be re-pointed to the 15th of the
# NEXT month (the actual mid-point "as-of" indexing date used by the
# naive forecaster). This recomputes each date in-place via the Excel
# object model rather than hard-coding serial numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 73
$col = "A"

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Range("$col$row")
    $serial = $cell.Value2

    if ($serial -eq $null) { continue }

    # Convert the Excel serial date (OLE Automation date) to a .NET DateTime
    $current = [DateTime]::FromOADate($serial)

    # Move to the same day-of-month in the following month, then pin the
    # day to the 15th (the corrected "quarterly index" anchor date).
    $nextMonth = $current.AddMonths(1)
    $fixed = Get-Date -Year $nextMonth.Year -Month $nextMonth.Month -Day 15 -Hour 0 -Minute 0 -Second 0

    $cell.Value2 = $fixed.ToOADate()
}
